$wb = $excel.ActiveWorkbook

$itemSheet = $wb.Worksheets.Item("Item")

# Add new "Order" column (column E) to the Item sheet
$itemSheet.Range("E1").Value = "int"
$itemSheet.Range("E2").Value = "Order"
$itemSheet.Range("E3").Value = 1
$itemSheet.Range("E4").Value = 1
$itemSheet.Range("E5").Value = 2
$itemSheet.Range("E6").Value = 3
$itemSheet.Range("E7").Value = 4
$itemSheet.Range("E8").Value = 5

# Match style of header row (E1) and second header row (E2) to neighboring cells
$itemSheet.Range("A1").Copy() | Out-Null
$itemSheet.Range("E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$itemSheet.Range("A2").Copy() | Out-Null
$itemSheet.Range("E2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Switch active sheet to Item, with selection at F5
$itemSheet.Activate() | Out-Null
$itemSheet.Range("F5").Select() | Out-Null
